# Apply cryptos list price/volume updates (GitHub Actions refresh).
# D-column values are decimal-looking strings that must stay as TEXT
# (matching the original inlineStr cells), so NumberFormat is forced to
# "@" before the assignment (otherwise COM auto-converts them to floats,
# e.g. "1.001" -> 1.0009999999999999) and the style is reset to "Normal"
# afterwards so no residual cell-level formatting is introduced.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "26.331.16"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -4.48%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.762.52"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -4.01%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.001"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.01%  "

$ws.Range("E5").Value = "  +0.02%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "304.25"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -2.57%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4259"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -0.47%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3614"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -1.09%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.07027"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -3.34%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.8301"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -4.00%  "

$ws.Range("E11").Value = "  -2.64%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.725.06"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -5.53%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "5.235"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -4.27%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.382"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -2.17%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.06782"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -2.75%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "1.002"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +0.00%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "79.07"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -2.00%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.000008652"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -2.70%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "1.001"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +0.06%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "14.96"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -2.98%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "26.330.92"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -4.33%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.990"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -3.25%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "11.09"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +1.83%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "1.977.24"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -2.87%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "1.906"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -4.25%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "152.00"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -1.86%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "18.10"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -4.08%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "114.61"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +0.24%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "5.006"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -3.07%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.668"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -8.33%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.08882"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +0.29%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.7213"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -3.81%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.112"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -1.96%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "4.298"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -5.55%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.001"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +0.01%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "2.713"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -9.45%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.067"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -2.83%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.05090"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -4.39%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.01878"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -3.11%  "

$ws.Range("E40").Value = "  -3.01%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.4888"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -3.71%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "6.173"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -4.71%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "2.485"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -11.20%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "7.957"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -4.60%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "104.53"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -0.99%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "1.001"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +0.10%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "9.982"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -4.57%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.06182"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -4.55%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.4454"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -5.13%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.705"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -2.05%  "
